# "updated koper libs files"
#
# Add new KOPER_IND_220..225 counterparty indicator rows to the
# "r AnalysisUnit_Variable" sheet (the 2nd worksheet in the workbook).
#
# Each existing data row in that sheet has:
#   A = Action            ("CREATE/MODIFY")
#   B = Id                 (e.g. "COUNTERPARTY_KOPER_IND_219")
#   C = Name               (same value as Id)
#   D = (unused / empty)
#   E = AnalysisUnit ref   ("COUNTERPARTY_KOPER")
#   F = Variable           (e.g. "KOPER_IND_219")

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$templateRow = 133
$lastRow = 133

# New counterparty indicator rows to append (in authoring order).
$newRows = @(
    @{ Id = "COUNTERPARTY_KOPER_IND_220"; Variable = "KOPER_IND_220"; VariableFirst = $true  },
    @{ Id = "COUNTERPARTY_KOPER_IND_221"; Variable = "KOPER_IND_221"; VariableFirst = $false },
    @{ Id = "COUNTERPARTY_KOPER_IND_222"; Variable = "KOPER_IND_222"; VariableFirst = $false },
    @{ Id = "COUNTERPARTY_KOPER_IND_223"; Variable = "KOPER_IND_223"; VariableFirst = $false },
    @{ Id = "COUNTERPARTY_KOPER_IND_224"; Variable = "KOPER_IND_224"; VariableFirst = $false },
    @{ Id = "COUNTERPARTY_KOPER_IND_225"; Variable = "KOPER_IND_225"; VariableFirst = $false }
)

$r = $lastRow
foreach ($row in $newRows) {
    $r = $r + 1

    if ($row.VariableFirst) {
        $ws2.Cells.Item($r, 6).Value = $row.Variable
        $ws2.Cells.Item($r, 2).Value = $row.Id
        $ws2.Cells.Item($r, 3).Value = $row.Id
    } else {
        $ws2.Cells.Item($r, 2).Value = $row.Id
        $ws2.Cells.Item($r, 3).Value = $row.Id
        $ws2.Cells.Item($r, 6).Value = $row.Variable
    }

    $ws2.Cells.Item($r, 1).Value = "CREATE/MODIFY"
    $ws2.Cells.Item($r, 5).Value = "COUNTERPARTY_KOPER"

    # Copy formatting (cell style) from the template row so the new rows
    # match the look of the existing ones.
    $ws2.Cells.Item($templateRow, 1).Copy()
    $ws2.Cells.Item($r, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# Update the on-screen selection to reflect where the author ended up
# after entering the new data (sheet2 scrolled/selected at the new last
# row, sheet1 left with its own last-used selection).
[void]$ws2.Range("B135").Select()
[void]$ws1.Range("E29").Select()
